$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 800, pushing existing rows
# 800-888 down to 802-890 (mirrors the two brand-new price records that
# were added to the weekly "Papa" series for Macroferia Regional de Talca).
$ws.Rows(800).Insert()
$ws.Rows(800).Insert()

# New row 800: Asterix, 1a (cosecha)
$ws.Range("A800").Value = 5
$ws.Range("B800").Value = "Macroferia Regional de Talca"
$ws.Range("C800").Value = "Maule"
$ws.Range("D800").Value = 45142
$ws.Range("E800").Value = 7
$ws.Range("F800").Value = 100114001
$ws.Range("G800").Value = "Papa"
$ws.Range("H800").Value = "Asterix"
$ws.Range("I800").Value = "1a (cosecha)"
$ws.Range("J800").Value = 1500
$ws.Range("K800").Value = 18000
$ws.Range("L800").Value = 18000
$ws.Range("M800").Value = 18000
$ws.Range("N800").Value = "`$/saco 25 kilos"
$ws.Range("O800").Value = "Región del Maule"
$ws.Range("P800").Value = 720
$ws.Range("Q800").Value = 25
$ws.Range("R800").Value = "Hortaliza"

# New row 801: Rosara, 1a (cosecha)
$ws.Range("A801").Value = 5
$ws.Range("B801").Value = "Macroferia Regional de Talca"
$ws.Range("C801").Value = "Maule"
$ws.Range("D801").Value = 45142
$ws.Range("E801").Value = 7
$ws.Range("F801").Value = 100114001
$ws.Range("G801").Value = "Papa"
$ws.Range("H801").Value = "Rosara"
$ws.Range("I801").Value = "1a (cosecha)"
$ws.Range("J801").Value = 1500
$ws.Range("K801").Value = 17000
$ws.Range("L801").Value = 17000
$ws.Range("M801").Value = 17000
$ws.Range("N801").Value = "`$/saco 25 kilos"
$ws.Range("O801").Value = "Región del Maule"
$ws.Range("P801").Value = 680
$ws.Range("Q801").Value = 25
$ws.Range("R801").Value = "Hortaliza"
